# Apply the edits described by the diff:
#  - F15, F16, F19 get value 5
#  - Active selection on the sheet moves to L14 (pane/topLeftCell follow automatically)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 5
$ws.Range("F19").Value = 5

$ws.Range("L14").Select()
